{"js": "// Replace specific paragraph bodies with redacted placeholder tokens.\n// Targets (by original, pre-edit text) are matched individually so that\n// duplicate paragraphs elsewhere in the document are left untouched:\n//   1) first \"P(x = 0) = 0.5^5 = 1/32 = 3.125%\"                     -> \"[CONTENT]\"\n//   2) \"[ANSWER]\"                                                    -> \"[CONTENT]\"\n//   3) \"The probability that she gets bumped is 1.40/1000, ...\"      -> \"[CALCULATED_VALUE]\"\n//   4) \"Answer: For independent travelers, the probability of ...\"   -> \"Answer: [CONTENT]\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst replacements = [\n  { match: \"P(x = 0) = 0.5^5 = 1/32 = 3.125%\", newText: \"[CONTENT]\", once: true },\n  { match: \"[ANSWER]\", newText: \"[CONTENT]\", once: true },\n  {\n    match:\n      \"The probability that she gets bumped is 1.40/1000, which equals 0.14%. This means that for every 1000 opportunities, there is a 1.4 chance of her getting bumped. For example, if she were to enter a crowd of 1000 people, she would statistically expect to get bumped about 1.4 times.\",\n    newText: \"[CALCULATED_VALUE]\",\n    once: true,\n  },\n  {\n    match:\n      \"Answer: For independent travelers, the probability of being bumped from a flight is calculated as 1 - (1 - 0.14)^10. This results in approximately 0.0139, or about 1.39%. This means that if a traveler books 10 flights, there is a 1.39% chance of being bumped from at least one of those flights.\",\n    newText: \"Answer: [CONTENT]\",\n    once: true,\n  },\n];\n\nfor (const rule of replacements) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    const para = paragraphs.items[i];\n    if (para.text === rule.match) {\n      para.insertText(rule.newText, Word.InsertLocation.replace);\n      if (rule.once) break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace specific paragraph bodies with redacted placeholder tokens.\n# Targets (by original, pre-edit text) are matched individually so that\n# duplicate paragraphs elsewhere in the document are left untouched:\n#   1) first \"P(x = 0) = 0.5^5 = 1/32 = 3.125%\"                     -> \"[CONTENT]\"\n#   2) \"[ANSWER]\"                                                    -> \"[CONTENT]\"\n#   3) \"The probability that she gets bumped is 1.40/1000, ...\"      -> \"[CALCULATED_VALUE]\"\n#   4) \"Answer: For independent travelers, the probability of ...\"   -> \"Answer: [CONTENT]\"\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Match = \"P(x = 0) = 0.5^5 = 1/32 = 3.125%\"; New = \"[CONTENT]\" },\n    @{ Match = \"[ANSWER]\"; New = \"[CONTENT]\" },\n    @{ Match = \"The probability that she gets bumped is 1.40/1000, which equals 0.14%. This means that for every 1000 opportunities, there is a 1.4 chance of her getting bumped. For example, if she were to enter a crowd of 1000 people, she would statistically expect to get bumped about 1.4 times.\"; New = \"[CALCULATED_VALUE]\" },\n    @{ Match = \"Answer: For independent travelers, the probability of being bumped from a flight is calculated as 1 - (1 - 0.14)^10. This results in approximately 0.0139, or about 1.39%. This means that if a traveler books 10 flights, there is a 1.39% chance of being bumped from at least one of those flights.\"; New = \"Answer: [CONTENT]\" }\n)\n\nforeach ($rule in $replacements) {\n    foreach ($p in $d.Paragraphs) {\n        $rng = $p.Range\n        $text = $rng.Text\n        if ($text.Length -gt 0 -and $text.Substring(0, $text.Length - 1) -eq $rule.Match) {\n            $rng.Text = $rule.New\n            break\n        }\n    }\n}\n"}
